$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The case results for the 380 kV line case were recomputed; update the
# pl_mw sheet values in range B2:O25 accordingly. Columns D,H,I,J,L,N
# are untouched (remain 0), so read the existing block first and only
# overwrite the cells whose values actually changed.
$range = $ws.Range("B2:O25")
$data = $range.Value2

$data[1,1] = 0.2865685268086793
$data[1,2] = 0.04897208787983232
$data[1,4] = 0.1637082793371789
$data[1,5] = 0.4443680307746263
$data[1,6] = 0.002443763844237843
$data[1,10] = 0.2588816154759286
$data[1,12] = 0.2190209090369351
$data[1,14] = 2.706722623918409
$data[2,1] = 0.2543081462343935
$data[2,2] = 0.04591971149869067
$data[2,4] = 0.1527014476912427
$data[2,5] = 0.387822817061874
$data[2,6] = 0.002446159046098492
$data[2,10] = 0.2262455101951701
$data[2,12] = 0.197160567764243
$data[2,14] = 2.734905478134507
$data[3,1] = 0.2345051496715769
$data[3,2] = 0.04403359983195543
$data[3,4] = 0.1460592999685133
$data[3,5] = 0.3531389305168915
$data[3,6] = 0.002447706873416658
$data[3,10] = 0.2061591229468434
$data[3,12] = 0.1838068273700202
$data[3,14] = 2.75396649619357
$data[4,1] = 0.2264369429234421
$data[4,2] = 0.04326203958245145
$data[4,4] = 0.1433815399569269
$data[4,5] = 0.3390132514313251
$data[4,6] = 0.002448357087439793
$data[4,10] = 0.1979621813880073
$data[4,12] = 0.1783822919790268
$data[4,14] = 2.762175173148051
$data[5,1] = 0.2250973368183509
$data[5,2] = 0.04313374557222005
$data[5,4] = 0.1429386431015871
$data[5,5] = 0.336668177824194
$data[5,6] = 0.002448466231829201
$data[5,10] = 0.1966003986439517
$data[5,12] = 0.1774825914512945
$data[5,14] = 2.763564846066259
$data[6,1] = 0.2343963316776012
$data[6,2] = 0.04402320619887234
$data[6,4] = 0.1460230697597567
$data[6,5] = 0.3529483938344953
$data[6,6] = 0.002447715563534969
$data[6,10] = 0.2060486224427649
$data[6,12] = 0.1837336004904699
$data[6,14] = 2.754075415804309
$data[7,1] = 0.2754443471152968
$data[7,2] = 0.04792213449003668
$data[7,4] = 0.1598888912119421
$data[7,5] = 0.4248636149813336
$data[7,6] = 0.002444573729695588
$data[7,10] = 0.2476388223577715
$data[7,12] = 0.2114692024292921
$data[7,14] = 2.716075145327636
$data[8,1] = 0.3559648970704359
$data[8,2] = 0.05547143410115041
$data[8,4] = 0.188012417002021
$data[8,5] = 0.5661985755041457
$data[8,6] = 0.002439022193594522
$data[8,10] = 0.3288044482970349
$data[8,12] = 0.2664070789901629
$data[8,14] = 2.655520537656031
$data[9,1] = 0.4151257780972344
$data[9,2] = 0.06095727190802336
$data[9,4] = 0.2092615246864753
$data[9,5] = 0.6702781546542269
$data[9,6] = 0.002435311352219402
$data[9,10] = 0.3881842087476457
$data[9,12] = 0.3071142568907419
$data[9,14] = 2.61957993469008
$data[10,1] = 0.4420377402954898
$data[10,2] = 0.06343941120111651
$data[10,4] = 0.2190597369055638
$data[10,5] = 0.7176906081379002
$data[10,6] = 0.002433702275365839
$data[10,10] = 0.4151403310059436
$data[10,12] = 0.3257101776984115
$data[10,14] = 2.605093791019215
$data[11,1] = 0.4522281795919127
$data[11,2] = 0.06437736648614134
$data[11,4] = 0.2227893081077994
$data[11,5] = 0.7356546913071611
$data[11,6] = 0.002433104260176894
$data[11,10] = 0.4253395248437641
$data[11,12] = 0.3327632893757624
$data[11,14] = 2.599876873644206
$data[12,1] = 0.4500335171306631
$data[12,2] = 0.06417544976144995
$data[12,4] = 0.2219852204973023
$data[12,5] = 0.7317853510981394
$data[12,6] = 0.002433232551456401
$data[12,10] = 0.4231433313080117
$data[12,12] = 0.3312437753392601
$data[12,14] = 2.600988469570041
$data[13,1] = 0.4428761271816768
$data[13,2] = 0.06351661725665281
$data[13,4] = 0.2193661848071287
$data[13,5] = 0.7191683204515869
$data[13,6] = 0.00243365284978266
$data[13,10] = 0.4159795978210923
$data[13,12] = 0.3262902160527759
$data[13,14] = 2.604659203435403
$data[14,1] = 0.4384919360317383
$data[14,2] = 0.06311280436537459
$data[14,4] = 0.2177644573480393
$data[14,5] = 0.7114413442032514
$data[14,6] = 0.00243391176575172
$data[14,10] = 0.4115904834184505
$data[14,12] = 0.3232574823702024
$data[14,14] = 2.606942643973611
$data[15,1] = 0.41336696728888
$data[15,2] = 0.06079478436561203
$data[15,4] = 0.2086238625162409
$data[15,5] = 0.6671810134426437
$data[15,6] = 0.002435418094208472
$data[15,10] = 0.3864213969766865
$data[15,12] = 0.3059005414744931
$data[15,14] = 2.62056419597576
$data[16,1] = 0.3979531348144008
$data[16,2] = 0.05936928694701749
$data[16,4] = 0.2030503506180352
$data[16,5] = 0.6400460337125793
$data[16,6] = 0.002436362374610064
$data[16,10] = 0.3709663078326173
$data[16,12] = 0.2952726180772274
$data[16,14] = 2.629398350601662
$data[17,1] = 0.3890874743619577
$data[17,2] = 0.05854812065138049
$data[17,4] = 0.1998570136699627
$data[17,5] = 0.6244449056556647
$data[17,6] = 0.002436912939147206
$data[17,10] = 0.362071702227837
$data[17,12] = 0.2891670601101097
$data[17,14] = 2.634654873619738
$data[18,1] = 0.3860857223055518
$data[18,2] = 0.05826987336494938
$data[18,4] = 0.1987779271254126
$data[18,5] = 0.619163680173358
$data[18,6] = 0.002437100630175984
$data[18,10] = 0.359059254099293
$data[18,12] = 0.2871010835466734
$data[18,14] = 2.636464736946138
$data[19,1] = 0.3995939701265172
$data[19,2] = 0.05952116413574515
$data[19,4] = 0.203642375448517
$data[19,5] = 0.642933953830422
$data[19,6] = 0.002436261084823754
$data[19,10] = 0.3726120747649304
$data[19,12] = 0.2964032186020304
$data[19,14] = 2.628439786419108
$data[20,1] = 0.4449784442295481
$data[20,2] = 0.06371018639120507
$data[20,4] = 0.2201349360889324
$data[20,5] = 0.7228739723491628
$data[20,6] = 0.002433529091482304
$data[20,10] = 0.4180839943650767
$data[20,12] = 0.3277448912715712
$data[20,14] = 2.60357372258153
$data[21,1] = 0.4746364872434015
$data[21,2] = 0.0664364079856341
$data[21,4] = 0.2310258049445508
$data[21,5] = 0.7751780083420101
$data[21,6] = 0.002431809459539285
$data[21,10] = 0.4477527452445429
$data[21,12] = 0.3482940118520546
$data[21,14] = 2.588888662872023
$data[22,1] = 0.4588078838820877
$data[22,2] = 0.0649824453449952
$data[22,4] = 0.2252028148064937
$data[22,5] = 0.7472568307830727
$data[22,6] = 0.002432721248689765
$data[22,10] = 0.4319226763321922
$data[22,12] = 0.3373205568535909
$data[22,14] = 2.596582798637741
$data[23,1] = 0.3988521610421856
$data[23,2] = 0.05945250552868231
$data[23,4] = 0.2033746869375577
$data[23,5] = 0.6416283278902171
$data[23,6] = 0.002436306854039146
$data[23,10] = 0.3718680523743387
$data[23,12] = 0.2958920598538484
$data[23,14] = 2.628872599961824
$data[24,1] = 0.3341805899470671
$data[24,2] = 0.05343966738962536
$data[24,4] = 0.1803025310304918
$data[24,5] = 0.5279251897347166
$data[24,6] = 0.002440459159533538
$data[24,10] = 0.3068904131806676
$data[24,12] = 0.2514851848127293
$data[24,14] = 2.670403505629537

$range.Value2 = $data

Write-Host "case with 380 kV done"
